# AFASpectraCalculations.xlsx edit:
# "Split of AFASpectra and reduction to max 2048 bytes"
#
# Rebuilds the worksheet: a new "which one / numSamplesSDSS" selector block,
# an INDEX()-based numSamples formula, a Name/Type/TypeSize/NumOf/Required/Total
# table with a Required flag (col F) and IF(...) Total column (col G), and a
# grand Total row at the bottom.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Start from a clean sheet - the old layout (B4:F19) is being restructured.
$ws.Cells.Clear()

# --- New header/selector block (rows 2-3), then the existing BOSS inputs ---
# Shared-string insertion order matters (matches the target sharedStrings.xml
# ordering), so write B3/B2 before the table header strings below.
$ws.Range("B3").Value = "numSamplesSDSS"
$ws.Range("G3").Value = 3900

$ws.Range("B2").Value = "which?"
$ws.Range("G2").Value = 1

$ws.Range("B4").Value = "numSamplesBOSS"
$ws.Range("G4").Value = 4700

$ws.Range("B5").Value = "reductionFactor"
$ws.Range("G5").Value = 8

$ws.Range("B6").Value = "numSamples"
$ws.Range("G6").Formula = "=ROUNDDOWN(INDEX(G3:G4,G2,1)/G5,0)"

# --- Table header row ---
$ws.Range("B8").Value = "Name"
$ws.Range("C8").Value = "Type"
$ws.Range("D8").Value = "TypeSize"
$ws.Range("E8").Value = "NumOf"
$ws.Range("F8").Value = "Required"
$ws.Range("G8").Value = "Total"

# --- Row 9: m_Amplitude (NumOf references numSamples via G6) ---
$ws.Range("B9").Value = "m_Amplitude"
$ws.Range("C9").Value = "float"
$ws.Range("D9").Value = 4
$ws.Range("E9").Formula = "=G6"
$ws.Range("F9").Value = $true
$ws.Range("G9").Formula = "=IF(F9,D9*E9,0)"

# --- Rows 10-19: remaining struct fields ---
$ws.Range("B10").Value = "m_Min"
$ws.Range("C10").Value = "float"
$ws.Range("D10").Value = 4
$ws.Range("E10").Value = 1
$ws.Range("F10").Value = $true

$ws.Range("B11").Value = "m_Max"
$ws.Range("C11").Value = "float"
$ws.Range("D11").Value = 4
$ws.Range("E11").Value = 1
$ws.Range("F11").Value = $true

$ws.Range("B12").Value = "m_Index"
$ws.Range("C12").Value = "sint32_t"
$ws.Range("D12").Value = 4
$ws.Range("E12").Value = 1
$ws.Range("F12").Value = $false

$ws.Range("B13").Value = "m_SamplesRead"
$ws.Range("C13").Value = "sint16_t"
$ws.Range("D13").Value = 2
$ws.Range("E13").Value = 1
$ws.Range("F13").Value = $false

$ws.Range("B14").Value = "m_SpecObjID"
$ws.Range("C14").Value = "sint64_t"
$ws.Range("D14").Value = 8
$ws.Range("E14").Value = 1
$ws.Range("F14").Value = $true

$ws.Range("B15").Value = "m_Type"
$ws.Range("C15").Value = "sint32_t"
$ws.Range("D15").Value = 4
$ws.Range("E15").Value = 1
$ws.Range("F15").Value = $false

$ws.Range("B16").Value = "m_version"
$ws.Range("C16").Value = "SpectraVersion"
$ws.Range("D16").Value = 4
$ws.Range("E16").Value = 1
$ws.Range("F16").Value = $false

$ws.Range("B17").Value = "m_Z"
$ws.Range("C17").Value = "double"
$ws.Range("D17").Value = 8
$ws.Range("E17").Value = 1
$ws.Range("F17").Value = $false

$ws.Range("B18").Value = "m_flux"
$ws.Range("C18").Value = "float"
$ws.Range("D18").Value = 4
$ws.Range("E18").Value = 1
$ws.Range("F18").Value = $true

$ws.Range("B19").Value = "m_status"
$ws.Range("C19").Value = "char"
$ws.Range("D19").Value = 1
$ws.Range("E19").Value = 1
$ws.Range("F19").Value = $false

# Shared formula across the whole Total column for rows 10-19.
$ws.Range("G10:G19").Formula = "=IF(F10,D10*E10,0)"

# --- Grand total row ---
$ws.Range("G21").Formula = "=SUM(G9:G19)"

# --- Column widths (best-fit approximation; the engine quantizes
#     ColumnWidth to 1/6-character steps with a fixed padding offset, so
#     these inputs are pre-compensated to land as close as possible to the
#     real Excel "bestFit" widths recorded in the target file) ---
$ws.Columns("B").ColumnWidth = 14.111979166666666
$ws.Columns("C").ColumnWidth = 11.744791666666666
$ws.Columns("D").ColumnWidth = 6.690104166666667
$ws.Columns("E").ColumnWidth = 5.584635416666667
$ws.Columns("F").ColumnWidth = 6.901041666666667
$ws.Columns("G").ColumnWidth = 3.8502604166666665

# --- Selection / view state ---
[void]$ws.Range("A3").Select()
